$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B57").Value = "09ecd661d724a2aad45dafbb1c481fb7"
$ws.Range("B154").Value = "0164192226833e8b2508d9634b0ba903"
$ws.Range("B216").Value = "86df63785cf20188f1796a5abd02a6b5"
$ws.Range("B273").Value = "0e4158b3be5756e9866cace2776c8db4"
$ws.Range("B299").Value = "ca06a29ddf84c1012ce23445464311d1"
$ws.Range("B330").Value = "02d08555a89aca4227289c60c19d9b82"
$ws.Range("B350").Value = "205045de71ccf4d8ebb7043be63d7d1e"
$ws.Range("B387").Value = "a6232cc167e916c3bae255a3aa1b496e"
$ws.Range("B424").Value = "c3d15ba386f49a4a89cff768392ffa95"
$ws.Range("B552").Value = "b87b0ff9b1bd0957496b465abc3e1606"
$ws.Range("B601").Value = "1aea55cc5703b249fea06d459a96cf71"
$ws.Range("B655").Value = "3379e70f93178a55f709d366d220e3ba"
$ws.Range("B712").Value = "f5c07954d5e36d9a67fc8c20c5548bcb"
$ws.Range("B731").Value = "a561d1bf4aefcf39e61e1863b8147b44"
$ws.Range("B740").Value = "d4374f0fa39c6f7edfbd28cca214f2b8"
$ws.Range("B802").Value = "11e6135d92906710ca6283d07f1d1add"
$ws.Range("B811").Value = "5f1e48ea2ee37ac4a0cd6534daf28e1d"
$ws.Range("B839").Value = "838e687b650fda7a6da60c9e4c56a4be"
$ws.Range("B846").Value = "da70563953f6e5c1d4a1aab0bbe1d7e0"
$ws.Range("B848").Value = "661c7a2286dd8390bd5f9d2ff11d671b"
$ws.Range("B874").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
$ws.Range("B911").Value = "cba30d7950a13a0c0967661dd8f1ded3"
$ws.Range("B951").Value = "3f574683856d8cc29639b08f7ab41e07"
